$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Through 2022-03-12"

$ws.Range("I1").Value = "2022 (through 03-12)"
$ws.Range("I4").Value = 56
$ws.Range("I14").Value = 356
